$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(0.091144958878133092, 0.82415920728107162, 0.53398025938135052, 0.95475415217332626),
    @(0.20103665446257424,  0.71552014005783737, 0.36889142831119193, 0.90145082795890596),
    @(0.1094911381302848,   0.84357644622325689, 0.39372854593393602, 0.94543545454701361),
    @(0.050214227820487532, 0.79888161464938912, 0.73672887591025371, 0.97501990047007026)
)

for ($r = 0; $r -lt 4; $r++) {
    for ($c = 0; $c -lt 4; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $values[$r][$c]
    }
}

# Target stored widths are 13.7109375 (col A) and 12.7109375 (cols B-D).
# The host only stores column width at whole-pixel granularity (width*6 must
# be an integer), so those exact values are unreachable; 12.76/11.76 are the
# smallest ColumnWidth inputs that round to the nearest achievable stored
# widths (13.666666... / 12.666666...).
$ws.Columns.Item(1).ColumnWidth = 12.76
$ws.Columns.Item(2).ColumnWidth = 11.76
$ws.Columns.Item(3).ColumnWidth = 11.76
$ws.Columns.Item(4).ColumnWidth = 11.76
